# Auto-generated Excel COM-interop script
# Applies the scraped-schedule refresh described in the commit diff
# across all three sheets (LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")

# Row 2
$ws.Cells.Item(2,1).Value = "Última actualización: 21:56:33"

# Row 3
$ws.Cells.Item(3,1).Value = "Total filas: 564"

# Row 92
$ws.Cells.Item(92,3).Value = "23_HERNANDEZ"

# Row 93
$ws.Cells.Item(93,3).Value = "215A_EL PATO"

# Row 108
$ws.Cells.Item(108,1).Value = "07:46:15"
$ws.Cells.Item(108,3).Value = "16_SANTA ANA"
$ws.Cells.Item(108,4).Value = 96

# Row 109
$ws.Cells.Item(109,1).Value = "08:21:27"
$ws.Cells.Item(109,3).Value = "17_ROMERO"
$ws.Cells.Item(109,4).Value = 61

# Row 110
$ws.Cells.Item(110,1).Value = "07:59:05"
$ws.Cells.Item(110,3).Value = "16_SANTA ANA"
$ws.Cells.Item(110,4).Value = 84

# Row 111
$ws.Cells.Item(111,1).Value = "08:21:27"
$ws.Cells.Item(111,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(111,4).Value = 62

# Row 112
$ws.Cells.Item(112,1).Value = "07:46:15"
$ws.Cells.Item(112,3).Value = "17_ROMERO"
$ws.Cells.Item(112,4).Value = 97

# Row 206
$ws.Cells.Item(206,3).Value = "10_OLMOS"

# Row 207
$ws.Cells.Item(207,3).Value = "84_COLONIA URQUIZA-ESC 49"

# Row 214
$ws.Cells.Item(214,1).Value = "11:51:05"
$ws.Cells.Item(214,3).Value = "14_ABASTO"
$ws.Cells.Item(214,4).Value = 30

# Row 215
$ws.Cells.Item(215,1).Value = "12:16:51"
$ws.Cells.Item(215,3).Value = "16_SANTA ANA"
$ws.Cells.Item(215,4).Value = 5

# Row 216
$ws.Cells.Item(216,1).Value = "10:28:12"
$ws.Cells.Item(216,3).Value = "215A_EL PATO"
$ws.Cells.Item(216,4).Value = 113

# Row 217
$ws.Cells.Item(217,1).Value = "10:28:12"
$ws.Cells.Item(217,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(217,4).Value = 113

# Row 224
$ws.Cells.Item(224,1).Value = "10:57:58"
$ws.Cells.Item(224,3).Value = "17_179 Y 38"
$ws.Cells.Item(224,4).Value = 100

# Row 225
$ws.Cells.Item(225,3).Value = "23_HERNANDEZ"

# Row 226
$ws.Cells.Item(226,1).Value = "11:51:05"
$ws.Cells.Item(226,3).Value = "27_EL RETIRO"
$ws.Cells.Item(226,4).Value = 46

# Row 230
$ws.Cells.Item(230,3).Value = "16_SANTA ANA"

# Row 231
$ws.Cells.Item(231,3).Value = "10_OLMOS"

# Row 251
$ws.Cells.Item(251,1).Value = "12:44:21"
$ws.Cells.Item(251,3).Value = "10_OLMOS"
$ws.Cells.Item(251,4).Value = 37

# Row 252
$ws.Cells.Item(252,1).Value = "11:51:05"
$ws.Cells.Item(252,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(252,4).Value = 90

# Row 322
$ws.Cells.Item(322,1).Value = "14:16:51"
$ws.Cells.Item(322,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(322,4).Value = 97

# Row 323
$ws.Cells.Item(323,1).Value = "13:55:06"
$ws.Cells.Item(323,3).Value = "14_ABASTO"
$ws.Cells.Item(323,4).Value = 118

# Row 409
$ws.Cells.Item(409,3).Value = "10_OLMOS"

# Row 410
$ws.Cells.Item(410,3).Value = "15_ABASTO"

# Row 437
$ws.Cells.Item(437,3).Value = "23_HERNANDEZ"

# Row 438
$ws.Cells.Item(438,3).Value = "11_ETCHEVERRY"

# Row 463
$ws.Cells.Item(463,3).Value = "16_SANTA ANA"

# Row 464
$ws.Cells.Item(464,3).Value = "11X44_ETCHEVERRY"

# Row 480
$ws.Cells.Item(480,1).Value = "18:12:36"
$ws.Cells.Item(480,3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(480,4).Value = 118

# Row 482
$ws.Cells.Item(482,1).Value = "20:10:57"
$ws.Cells.Item(482,3).Value = "10_OLMOS"
$ws.Cells.Item(482,4).Value = 0

# Row 490
$ws.Cells.Item(490,3).Value = "23_HERNANDEZ"

# Row 491
$ws.Cells.Item(491,3).Value = "17_ROMERO"

# Row 506
$ws.Cells.Item(506,3).Value = "16_SANTA ANA"

# Row 507
$ws.Cells.Item(507,1).Value = "20:40:23"
$ws.Cells.Item(507,3).Value = "17_ROMERO"
$ws.Cells.Item(507,4).Value = 1

# Row 508
$ws.Cells.Item(508,3).Value = "16_P MOR-167 Y 521"

# Row 509
$ws.Cells.Item(509,1).Value = "18:47:49"
$ws.Cells.Item(509,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(509,4).Value = 114

# Row 516
$ws.Cells.Item(516,1).Value = "19:12:20"
$ws.Cells.Item(516,3).Value = "15_ABASTO"
$ws.Cells.Item(516,4).Value = 100

# Row 517
$ws.Cells.Item(517,1).Value = "19:48:17"
$ws.Cells.Item(517,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(517,4).Value = 64

# Row 533
$ws.Cells.Item(533,1).Value = "19:35:03"
$ws.Cells.Item(533,3).Value = "10_OLMOS"
$ws.Cells.Item(533,4).Value = 108

# Row 534
$ws.Cells.Item(534,1).Value = "20:29:22"
$ws.Cells.Item(534,3).Value = "15_ABASTO"
$ws.Cells.Item(534,4).Value = 54

# Row 552
$ws.Cells.Item(552,1).Value = "21:56:33"
$ws.Cells.Item(552,2).Value = "22:23"
$ws.Cells.Item(552,4).Value = 27

# Row 553
$ws.Cells.Item(553,1).Value = "20:40:23"
$ws.Cells.Item(553,2).Value = "22:25"
$ws.Cells.Item(553,4).Value = 105

# Row 554
$ws.Cells.Item(554,1).Value = "20:55:06"
$ws.Cells.Item(554,2).Value = "22:26"
$ws.Cells.Item(554,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(554,4).Value = 91

# Row 555
$ws.Cells.Item(555,1).Value = "21:56:33"
$ws.Cells.Item(555,2).Value = "22:26"
$ws.Cells.Item(555,3).Value = "16_SANTA ANA"
$ws.Cells.Item(555,4).Value = 30

# Row 556
$ws.Cells.Item(556,1).Value = "20:48:29"
$ws.Cells.Item(556,2).Value = "22:27"
$ws.Cells.Item(556,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(556,4).Value = 99

# Row 557
$ws.Cells.Item(557,1).Value = "21:56:33"
$ws.Cells.Item(557,2).Value = "22:27"
$ws.Cells.Item(557,3).Value = "16_SANTA ANA"
$ws.Cells.Item(557,4).Value = 31

# Row 558
$ws.Cells.Item(558,1).Value = "20:29:22"
$ws.Cells.Item(558,2).Value = "22:28"
$ws.Cells.Item(558,3).Value = "84_COLONIA URQUIZA-ESC 49"

# Row 559
$ws.Cells.Item(559,2).Value = "22:29"
$ws.Cells.Item(559,3).Value = "10_OLMOS"
$ws.Cells.Item(559,4).Value = 94

# Row 560
$ws.Cells.Item(560,1).Value = "20:40:23"
$ws.Cells.Item(560,2).Value = "22:30"
$ws.Cells.Item(560,3).Value = "10_OLMOS"
$ws.Cells.Item(560,4).Value = 110

# New rows appended at the end of LP1912
$ws.Cells.Item(561,1).Value = "21:56:33"
$ws.Cells.Item(561,2).Value = "22:36"
$ws.Cells.Item(561,3).Value = "10_OLMOS"
$ws.Cells.Item(561,4).Value = 40
$ws.Cells.Item(561,5).Value = "LP1912"

$ws.Cells.Item(562,1).Value = "20:40:23"
$ws.Cells.Item(562,2).Value = "22:39"
$ws.Cells.Item(562,3).Value = "215A_EL PATO"
$ws.Cells.Item(562,4).Value = 119
$ws.Cells.Item(562,5).Value = "LP1912"

$ws.Cells.Item(563,1).Value = "20:55:06"
$ws.Cells.Item(563,2).Value = "22:41"
$ws.Cells.Item(563,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(563,4).Value = 106
$ws.Cells.Item(563,5).Value = "LP1912"

$ws.Cells.Item(564,1).Value = "20:55:06"
$ws.Cells.Item(564,2).Value = "22:48"
$ws.Cells.Item(564,3).Value = "14_ABASTO"
$ws.Cells.Item(564,4).Value = 113
$ws.Cells.Item(564,5).Value = "LP1912"

$ws.Cells.Item(565,1).Value = "21:56:33"
$ws.Cells.Item(565,2).Value = "22:52"
$ws.Cells.Item(565,3).Value = "10_OLMOS"
$ws.Cells.Item(565,4).Value = 56
$ws.Cells.Item(565,5).Value = "LP1912"

$ws.Cells.Item(566,1).Value = "21:56:33"
$ws.Cells.Item(566,2).Value = "23:07"
$ws.Cells.Item(566,3).Value = "14X44_ABASTO"
$ws.Cells.Item(566,4).Value = 71
$ws.Cells.Item(566,5).Value = "LP1912"

$ws.Cells.Item(567,1).Value = "21:56:33"
$ws.Cells.Item(567,2).Value = "23:08"
$ws.Cells.Item(567,3).Value = "17_ROMERO"
$ws.Cells.Item(567,4).Value = 72
$ws.Cells.Item(567,5).Value = "LP1912"

$ws.Cells.Item(568,1).Value = "21:56:33"
$ws.Cells.Item(568,2).Value = "23:12"
$ws.Cells.Item(568,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(568,4).Value = 76
$ws.Cells.Item(568,5).Value = "LP1912"

$ws.Cells.Item(569,1).Value = "21:56:33"
$ws.Cells.Item(569,2).Value = "23:38"
$ws.Cells.Item(569,3).Value = "215_ALUAR"
$ws.Cells.Item(569,4).Value = 102
$ws.Cells.Item(569,5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")

# Row 2
$ws.Cells.Item(2,1).Value = "Última actualización: 21:56:33"

# Row 3
$ws.Cells.Item(3,1).Value = "Total filas: 52"

# New rows appended at the end of LP1912-215
$ws.Cells.Item(57,1).Value = "21:56:33"
$ws.Cells.Item(57,2).Value = "23:38"
$ws.Cells.Item(57,3).Value = "215_ALUAR"
$ws.Cells.Item(57,4).Value = 102
$ws.Cells.Item(57,5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")

# Row 2
$ws.Cells.Item(2,1).Value = "Última actualización: 21:56:33"

# Row 3
$ws.Cells.Item(3,1).Value = "Total filas: 71"

# New rows appended at the end of 6203-6173
$ws.Cells.Item(76,1).Value = "21:56:33"
$ws.Cells.Item(76,2).Value = "23:08"
$ws.Cells.Item(76,3).Value = "215A_LA PLATA"
$ws.Cells.Item(76,4).Value = 72
$ws.Cells.Item(76,5).Value = "L6173"
